# Auto-generated edit script: apply scheduled-runner value updates
# to the Bahamut_Profits crafting-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 865.8333
$ws.Range("I103").Value = 998.75
$ws.Range("J103").Value = 600
$ws.Range("K103").Value = 2996.25
$ws.Range("L103").Value = 1800
$ws.Range("M103").Value = -2410.25
$ws.Range("N103").Value = -2972

$ws.Range("H116").Value = 5615
$ws.Range("I116").Value = 5916.6665
$ws.Range("K116").Value = 5916.6665
$ws.Range("M116").Value = -2474.6665

$ws.Range("H132").Value = 2122.1714
$ws.Range("I132").Value = 2476.0908
$ws.Range("K132").Value = 7428.2724
$ws.Range("M132").Value = -4898.2724

$ws.Range("H141").Value = 3511.875
$ws.Range("I141").Value = 3580.9092
$ws.Range("J141").Value = 3360
$ws.Range("K141").Value = 10742.7276
$ws.Range("L141").Value = 10080
$ws.Range("M141").Value = -5562.7276
$ws.Range("N141").Value = -20440

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 14367.1
$ws.Range("I28").Value = 5733.875
$ws.Range("J28").Value = 48900
$ws.Range("K28").Value = 5733.875
$ws.Range("L28").Value = 48900
$ws.Range("M28").Value = -5541.875
$ws.Range("N28").Value = -49284

$ws.Range("H99").Value = 14367.1
$ws.Range("I99").Value = 5733.875
$ws.Range("J99").Value = 48900
$ws.Range("K99").Value = 5733.875
$ws.Range("L99").Value = 48900
$ws.Range("M99").Value = -2738.875
$ws.Range("N99").Value = -54890

$ws.Range("H138").Value = 84940
$ws.Range("J138").Value = 84940
$ws.Range("L138").Value = 84940
$ws.Range("N138").Value = -95220

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 21727.54
$ws.Range("I134").Value = 1713.4722
$ws.Range("J134").Value = 73192.28999999999
$ws.Range("K134").Value = 5140.4166
$ws.Range("L134").Value = 219576.87
$ws.Range("M134").Value = -2605.4166
$ws.Range("N134").Value = -224646.87

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 334653.34
$ws.Range("J4").Value = 334653.34
$ws.Range("L4").Value = 334653.34
$ws.Range("N4").Value = -334877.34

$ws.Range("H99").Value = 1816
$ws.Range("I99").Value = 1615.92
$ws.Range("J99").Value = 2649.6667
$ws.Range("K99").Value = 1615.92
$ws.Range("L99").Value = 2649.6667
$ws.Range("M99").Value = -117.9200000000001
$ws.Range("N99").Value = -5645.6667

$ws.Range("H126").Value = 1816
$ws.Range("I126").Value = 1615.92
$ws.Range("J126").Value = 2649.6667
$ws.Range("K126").Value = 4847.76
$ws.Range("L126").Value = 7949.000100000001
$ws.Range("M126").Value = -2377.76
$ws.Range("N126").Value = -12889.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4683.846
$ws.Range("I3").Value = 1441.4286
$ws.Range("J3").Value = 8466.666999999999
$ws.Range("K3").Value = 4324.2858
$ws.Range("L3").Value = 25400.001
$ws.Range("M3").Value = -4212.2858
$ws.Range("N3").Value = -25624.001

$ws.Range("H4").Value = 1861
$ws.Range("I4").Value = 116
$ws.Range("J4").Value = 2234.9285
$ws.Range("K4").Value = 348
$ws.Range("L4").Value = 6704.7855
$ws.Range("M4").Value = -236
$ws.Range("N4").Value = -6928.7855

$ws.Range("H107").Value = 707954.25
$ws.Range("I107").Value = 1833.3334
$ws.Range("J107").Value = 972749.6
$ws.Range("K107").Value = 5500.0002
$ws.Range("L107").Value = 2918248.8
$ws.Range("M107").Value = -3580.0002
$ws.Range("N107").Value = -2922088.8

$ws.Range("H136").Value = 65603.125
$ws.Range("I136").Value = 144212.86
$ws.Range("K136").Value = 432638.58
$ws.Range("M136").Value = -427538.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5087.375
$ws.Range("I5").Value = 2181.6
$ws.Range("J5").Value = 6408.1816
$ws.Range("K5").Value = 2181.6
$ws.Range("L5").Value = 6408.1816
$ws.Range("M5").Value = -2069.6
$ws.Range("N5").Value = -6632.1816

$ws.Range("H42").Value = 25193.334
$ws.Range("J42").Value = 25193.334
$ws.Range("L42").Value = 25193.334
$ws.Range("N42").Value = -26163.334

$ws.Range("H70").Value = 5288
$ws.Range("I70").Value = 5599.8
$ws.Range("J70").Value = 5028.1665
$ws.Range("K70").Value = 5599.8
$ws.Range("L70").Value = 5028.1665
$ws.Range("M70").Value = -5329.8
$ws.Range("N70").Value = -5568.1665

$ws.Range("H73").Value = 5288
$ws.Range("I73").Value = 5599.8
$ws.Range("J73").Value = 5028.1665
$ws.Range("K73").Value = 5599.8
$ws.Range("L73").Value = 5028.1665
$ws.Range("M73").Value = -4663.8
$ws.Range("N73").Value = -6900.1665

$ws.Range("H93").Value = 38666.668
$ws.Range("J93").Value = 38666.668
$ws.Range("L93").Value = 38666.668
$ws.Range("N93").Value = -42410.668

$ws.Range("H97").Value = 2133.3333
$ws.Range("I97").Value = 2533.3333
$ws.Range("J97").Value = 1333.3334
$ws.Range("K97").Value = 2533.3333
$ws.Range("L97").Value = 1333.3334
$ws.Range("M97").Value = -2037.3333
$ws.Range("N97").Value = -2325.3334

$ws.Range("H98").Value = 24899.5
$ws.Range("J98").Value = 24899.5
$ws.Range("L98").Value = 24899.5
$ws.Range("N98").Value = -30889.5

$ws.Range("H100").Value = 24000
$ws.Range("J100").Value = 24000
$ws.Range("L100").Value = 24000
$ws.Range("N100").Value = -26164

$ws.Range("H113").Value = 7320
$ws.Range("I113").Value = 10200
$ws.Range("J113").Value = 2520
$ws.Range("K113").Value = 10200
$ws.Range("L113").Value = 2520
$ws.Range("M113").Value = -8030
$ws.Range("N113").Value = -6860

$ws.Range("H115").Value = 25193.334
$ws.Range("J115").Value = 25193.334
$ws.Range("L115").Value = 25193.334
$ws.Range("N115").Value = -27543.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2980
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H22").Value = 1582.5
$ws.Range("J22").Value = 1652.0588
$ws.Range("L22").Value = 1652.0588
$ws.Range("N22").Value = -2242.0588

$ws.Range("H27").Value = 1582.5
$ws.Range("J27").Value = 1652.0588
$ws.Range("L27").Value = 1652.0588
$ws.Range("N27").Value = -1866.0588

$ws.Range("H111").Value = 40000
$ws.Range("J111").Value = 40000
$ws.Range("L111").Value = 40000
$ws.Range("N111").Value = -48180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5385
$ws.Range("I2").Value = 4995
$ws.Range("J2").Value = 5450
$ws.Range("K2").Value = 4995
$ws.Range("L2").Value = 5450
$ws.Range("M2").Value = -4883
$ws.Range("N2").Value = -5674

$ws.Range("H100").Value = 2350
$ws.Range("I100").Value = 2350
$ws.Range("K100").Value = 4700
$ws.Range("M100").Value = -4159
